$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Unhide rows 2-49 (they were hidden in the source sheet; the new
#     commit reveals them). We unhide BEFORE touching any cell values so
#     the engine doesn't compute a stray auto row-height on rows that are
#     still flagged hidden while their content changes. ---
for ($r = 2; $r -le 49; $r++) {
    $ws.Rows.Item($r).Hidden = $false
}

# --- Row 3 : recompute of the preprocessing-threshold results for Sc1-3 ---
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = " 377075.00"
$ws.Range("D3").Style = "Normal"

$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = " 348407.82"
$ws.Range("E3").Style = "Normal"

$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 1

$ws.Range("I3").NumberFormat = "@"
$ws.Range("I3").Value = " 377075.00"
$ws.Range("I3").Style = "Normal"

$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 516
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = 0

$ws.Range("N3").NumberFormat = "@"
$ws.Range("N3").Value = " 0.66"
$ws.Range("N3").Style = "Normal"

$ws.Range("O3").Value = 216

# --- Row 4 : recompute of the preprocessing-threshold results for Sc1-4 ---
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = " 235565.00"
$ws.Range("D4").Style = "Normal"

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = " 233094.52"
$ws.Range("E4").Style = "Normal"

$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 1

$ws.Range("I4").NumberFormat = "@"
$ws.Range("I4").Value = " 235565.00"
$ws.Range("I4").Style = "Normal"

$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 444
$ws.Range("L4").Value = 73
$ws.Range("M4").Value = 0

$ws.Range("N4").NumberFormat = "@"
$ws.Range("N4").Value = " 0.61"
$ws.Range("N4").Style = "Normal"

$ws.Range("O4").Value = 204

# --- Move / record the active selection on the sheet ---
$ws.Range("K4").Select() | Out-Null
